# Update workbook data to match the target state:
# - Rows 2-14 values changed
# - Rows 15-35 added as new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for column A (Meta) and column B (Venda), rows 2-35
$data = @{
    2  = @(4000, 8000)
    3  = @(2000, 6000)
    4  = @(1000, 1000)
    5  = @(7777, 7777)
    6  = @(4000, 4000)
    7  = @(4000, 4000)
    8  = @(4575, 9985)
    9  = @(6000, 6000)
    10 = @(4000, 4000)
    11 = @(4000, 4000)
    12 = @(1780, 1898)
    13 = @(4000, 4000)
    14 = @(4000, 8000)
    15 = @(8040, 9000)
    16 = @(4000, 4000)
    17 = @(4000, 4000)
    18 = @(4000, 4000)
    19 = @(4000, 4000)
    20 = @(4000, 4000)
    21 = @(4000, 4000)
    22 = @(4000, 4000)
    23 = @(4000, 4000)
    24 = @(4000, 4000)
    25 = @(4000, 4000)
    26 = @(4000, 4000)
    27 = @(4000, 4000)
    28 = @(4000, 4000)
    29 = @(4000, 4000)
    30 = @(4000, 4000)
    31 = @(4000, 4000)
    32 = @(4000, 4000)
    33 = @(4000, 4000)
    34 = @(4000, 4000)
    35 = @(4000, 8000)
}

foreach ($row in $data.Keys | Sort-Object) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
}
